$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The naive forecaster data was regenerated: the oldest forecast row (dated
# 2007-11-14) is dropped, every remaining row shifts up by one, and the
# y_0_forecast / y_1_forecast columns are recomputed with corrected values,
# with one new forecast row appended at the end (dated 2025-08-22).
$ws.Rows.Item(2).Delete()

# A:E values for rows 2-52 after the shift (date serials, years, forecasts)
$data = New-Object 'object[,]' 51,5
$data[0,0] = 39583
$data[0,1] = 2008
$data[0,3] = 2009
$data[0,4] = 6.992417114397731
$data[1,0] = 39765
$data[1,1] = 2008
$data[1,3] = 2009
$data[1,4] = 2.651903832953884
$data[2,0] = 39948
$data[2,1] = 2009
$data[2,2] = 1.052572528399653
$data[2,3] = 2010
$data[2,4] = 2.396905458966625
$data[3,0] = 40130
$data[3,1] = 2009
$data[3,2] = 0.1715429114845346
$data[3,3] = 2010
$data[3,4] = 0.5356365903280924
$data[4,0] = 40310
$data[4,1] = 2010
$data[4,2] = 0.475547144596522
$data[4,3] = 2011
$data[4,4] = 1.495774185788745
$data[5,0] = 40494
$data[5,1] = 2010
$data[5,2] = 0.8004663283405655
$data[5,3] = 2011
$data[5,4] = 1.805987045940682
$data[6,0] = 40676
$data[6,1] = 2011
$data[6,2] = 4.312104569761943
$data[6,3] = 2012
$data[6,4] = 1.976172106438545
$data[7,0] = 40862
$data[7,1] = 2011
$data[7,2] = 5.253783907501819
$data[7,3] = 2012
$data[7,4] = 4.382531137514767
$data[8,0] = 41044
$data[8,1] = 2012
$data[8,2] = 4.639864760432189
$data[8,3] = 2013
$data[8,4] = 4.074582884048139
$data[9,0] = 41228
$data[9,1] = 2012
$data[9,2] = 3.522405026196918
$data[9,3] = 2013
$data[9,4] = 1.988448192515935
$data[10,0] = 41409
$data[10,1] = 2013
$data[10,2] = 1.418316910291906
$data[10,3] = 2014
$data[10,4] = 3.247860853607465
$data[11,0] = 41592
$data[11,1] = 2013
$data[11,2] = 1.656063945467268
$data[11,3] = 2014
$data[11,4] = 3.166945525867848
$data[12,0] = 41774
$data[12,1] = 2014
$data[12,2] = 3.520945360626571
$data[12,3] = 2015
$data[12,4] = 2.547371915279606
$data[13,0] = 41957
$data[13,1] = 2014
$data[13,2] = 4.06235252733802
$data[13,3] = 2015
$data[13,4] = 3.684750195712683
$data[14,0] = 42137
$data[14,1] = 2015
$data[14,2] = 3.285568146716344
$data[14,3] = 2016
$data[14,4] = 3.38738696315446
$data[15,0] = 42321
$data[15,1] = 2015
$data[15,2] = 3.05427116350534
$data[15,3] = 2016
$data[15,4] = 2.9587404276884
$data[16,0] = 42503
$data[16,1] = 2016
$data[16,2] = 2.441258738366514
$data[16,3] = 2017
$data[16,4] = 3.091878630346012
$data[17,0] = 42689
$data[17,1] = 2016
$data[17,2] = 2.305809238174006
$data[17,3] = 2017
$data[17,4] = 2.661643377950096
$data[18,0] = 42867
$data[18,1] = 2017
$data[18,2] = 2.480871685520603
$data[18,3] = 2018
$data[18,4] = 2.296583397191387
$data[19,0] = 43053
$data[19,1] = 2017
$data[19,2] = 2.509111342826809
$data[19,3] = 2018
$data[19,4] = 2.433248629349549
$data[20,0] = 43145
$data[20,1] = 2018
$data[20,2] = 2.83347664679956
$data[20,3] = 2019
$data[20,4] = 2.560065157976177
$data[21,0] = 43235
$data[21,1] = 2018
$data[21,2] = 3.071095202329288
$data[21,3] = 2019
$data[21,4] = 2.781797072072023
$data[22,0] = 43326
$data[22,1] = 2018
$data[22,2] = 3.248721852957415
$data[22,3] = 2019
$data[22,4] = 3.078223990352669
$data[23,0] = 43418
$data[23,1] = 2018
$data[23,2] = 3.296731496509198
$data[23,3] = 2019
$data[23,4] = 3.296423324101938
$data[24,0] = 43510
$data[24,1] = 2019
$data[24,2] = 3.149270133134596
$data[24,3] = 2020
$data[24,4] = 3.009352983329028
$data[25,0] = 43600
$data[25,1] = 2019
$data[25,2] = 2.994116795316071
$data[25,3] = 2020
$data[25,4] = 2.939737488252936
$data[26,0] = 43691
$data[26,1] = 2019
$data[26,2] = 2.867378798220366
$data[26,3] = 2020
$data[26,4] = 2.769017518462746
$data[27,0] = 43783
$data[27,1] = 2019
$data[27,2] = 2.861315725866587
$data[27,3] = 2020
$data[27,4] = 2.763966172716947
$data[28,0] = 43875
$data[28,1] = 2020
$data[28,2] = 2.689501145820206
$data[28,3] = 2021
$data[28,4] = 2.967409274751098
$data[29,0] = 43966
$data[29,1] = 2020
$data[29,2] = 2.671604274379558
$data[29,3] = 2021
$data[29,4] = 2.997455747043043
$data[30,0] = 44068
$data[30,1] = 2020
$data[30,2] = 1.790319754067715
$data[30,3] = 2021
$data[30,4] = 1.627017245406992
$data[31,0] = 44159
$data[31,1] = 2020
$data[31,2] = 1.790319754067715
$data[31,3] = 2021
$data[31,4] = 1.856930494010856
$data[32,0] = 44251
$data[32,1] = 2021
$data[32,2] = 1.980033360076905
$data[32,3] = 2022
$data[32,4] = 2.063021041451907
$data[33,0] = 44341
$data[33,1] = 2021
$data[33,2] = 2.08524086077817
$data[33,3] = 2022
$data[33,4] = 2.197771900625956
$data[34,0] = 44432
$data[34,1] = 2021
$data[34,2] = 2.339531676162721
$data[34,3] = 2022
$data[34,4] = 3.195599391913406
$data[35,0] = 44525
$data[35,1] = 2021
$data[35,2] = 2.339531676162721
$data[35,3] = 2022
$data[35,4] = 4.270817433327112
$data[36,0] = 44617
$data[36,1] = 2022
$data[36,2] = 4.865769161659883
$data[36,3] = 2023
$data[36,4] = 3.291462037299842
$data[37,0] = 44706
$data[37,1] = 2022
$data[37,2] = 4.939003803830477
$data[37,3] = 2023
$data[37,4] = 3.440178795466697
$data[38,0] = 44798
$data[38,1] = 2022
$data[38,2] = 4.834496776263886
$data[38,3] = 2023
$data[38,4] = 3.06038938938058
$data[39,0] = 44890
$data[39,1] = 2022
$data[39,2] = 4.834496776263886
$data[39,3] = 2023
$data[39,4] = 3.604316462518464
$data[40,0] = 44981
$data[40,1] = 2023
$data[40,2] = 3.332544669973525
$data[40,3] = 2024
$data[40,4] = 3.789179157493971
$data[41,0] = 45071
$data[41,1] = 2023
$data[41,2] = 2.93530792557688
$data[41,3] = 2024
$data[41,4] = 3.215749572764803
$data[42,0] = 45163
$data[42,1] = 2023
$data[42,2] = 2.798216547494237
$data[42,3] = 2024
$data[42,4] = 2.721520966738655
$data[43,0] = 45254
$data[43,1] = 2023
$data[43,2] = 2.798216547494237
$data[43,3] = 2024
$data[43,4] = 2.098908173995873
$data[44,0] = 45345
$data[44,1] = 2024
$data[44,2] = 1.888626610265987
$data[44,3] = 2025
$data[44,4] = 2.949781091571957
$data[45,0] = 45436
$data[45,1] = 2024
$data[45,2] = 1.635353376270698
$data[45,3] = 2025
$data[45,4] = 2.328770194687713
$data[46,0] = 45534
$data[46,1] = 2024
$data[46,2] = 1.530879676868468
$data[46,3] = 2025
$data[46,4] = 1.713178787950698
$data[47,0] = 45618
$data[47,1] = 2024
$data[47,2] = 1.530879676868468
$data[47,3] = 2025
$data[47,4] = 1.874466487556892
$data[48,0] = 45713
$data[48,1] = 2025
$data[48,2] = 1.966591496003445
$data[48,3] = 2026
$data[48,4] = 2.005435469818684
$data[49,0] = 45800
$data[49,1] = 2025
$data[49,2] = 1.984020855913604
$data[49,3] = 2026
$data[49,4] = 1.88544721086894
$data[50,0] = 45891
$data[50,1] = 2025
$data[50,2] = 2.060859685319461
$data[50,3] = 2026
$data[50,4] = 2.141985433296578

$ws.Range("A2:E52").Value = $data
